$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append row 3 with the new trip record (values stored as text, matching
# the existing rows which are all text-typed).
$row = 3

# A3 is blank in the source data - use a leading quote to force an empty
# text value (keeps the cell type Text instead of an empty/null Number).
$ws.Range("A$row").Value = "'"
$ws.Range("B$row").Value = "احمد"
# C3 looks numeric ("222") - force text so it isn't coerced to a number.
$ws.Range("C$row").Value = "'222"
$ws.Range("D$row").Value = "ايتا"
$ws.Range("E$row").Value = "الرحلة 3"
$ws.Range("F$row").Value = "C3"
$ws.Range("G$row").Value = "WCK"
$ws.Range("H$row").Value = "٠١‏/٠٥‏/٢٠٢٥ ٠٥:١٣:٢٩ م"
